$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.376.60"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.803.40"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.577"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.74%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "35.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.04%  "

$ws.Range("E9").Value = "  +2.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("D12").Value = "2.065.18"
$ws.Range("E12").Value = "  +1.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.42%  "

$ws.Range("D14").Value = "1.809.90"
$ws.Range("E14").Value = "  +1.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.09%  "

$ws.Range("D17").Value = "34.387.65"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("D20").Value = "0.0₃0795"
$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("E24").Value = "  +4.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.118"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.47%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.42%  "

$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("D35").Value = "1.396.54"
$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.674"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "

$ws.Range("E37").Value = "  -3.82%  "

$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("E39").Value = "  -0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.962"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "82.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0501"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.12%  "

$ws.Range("D48").Value = "1.965.56"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("E51").Value = "  +1.03%  "
